# Apply edits to the "kmeans_train" worksheet (legend/lookup table in columns I:M)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kmeans_train")

# --- "낮/밤" (day/night) legend block (rows 3-5): add a 3rd option "All Time" with code 2 ---
$ws.Range("K4").Value = 2
$ws.Range("K5").Value = "All Time"

# --- "계절" (season) legend block (rows 11-13): collapse 4 seasons into 2 combined seasons ---
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("I13").Value = "봄/여름"
$ws.Range("J13").Value = "가을/겨울"
$ws.Range("K13").ClearContents()
$ws.Range("L13").ClearContents()

# --- "상황" (situation) legend block (row 16): fix code values ---
$ws.Range("J16").Value = 1
$ws.Range("L16").Value = 3

# --- Update the selected cell shown when the sheet is reopened ---
$ws.Range("F14").Select()

# --- Restore the workbook window position ---
$excel.ActiveWindow.Left = 1635
$excel.ActiveWindow.Top = 2280
